$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-8
# from 2023-09-06 (45175) to 2023-09-14 (45183)
$newDate = 45183
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
